$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'256.55"
$ws.Range("E2").Value = "'-0.95%"
$ws.Range("D3").Value = "'27.27"
$ws.Range("E3").Value = "'-0.67%"
$ws.Range("D4").Value = "'4.668"
$ws.Range("E4").Value = "'-10.42%"
$ws.Range("E5").Value = "'-0.87%"
$ws.Range("D6").Value = "'6.630"
$ws.Range("E6").Value = "'-1.21%"
$ws.Range("D7").Value = "'0.8583"
$ws.Range("E7").Value = "'-1.49%"
$ws.Range("D8").Value = "'0.9400"
$ws.Range("E8").Value = "'-6.24%"
$ws.Range("D9").Value = "'0.1404"
$ws.Range("E9").Value = "'-1.19%"
$ws.Range("D10").Value = "'0.04106"
$ws.Range("E10").Value = "'15.23%"
$ws.Range("D11").Value = "'0.07103"
$ws.Range("E11").Value = "'-1.02%"
$ws.Range("D12").Value = "'0.03143"
$ws.Range("E12").Value = "'-0.25%"
$ws.Range("D13").Value = "'0.09145"
$ws.Range("E13").Value = "'-0.91%"
$ws.Range("D14").Value = "'0.001526"
$ws.Range("E14").Value = "'-0.93%"
$ws.Range("D15").Value = "'0.0006030"
$ws.Range("E15").Value = "'-0.60%"
$ws.Range("D16").Value = "'0.006213"
$ws.Range("E16").Value = "'6.84%"
$ws.Range("D17").Value = "'3.519"
$ws.Range("E17").Value = "'0.33%"
$ws.Range("E18").Value = "'-2.01%"
$ws.Range("E19").Value = "'0.36%"
$ws.Range("E20").Value = "'-2.92%"
$ws.Range("E21").Value = "'-0.48%"
$ws.Range("D22").Value = "'3.817"
$ws.Range("E22").Value = "'8.53%"
$ws.Range("D23").Value = "'0.04248"
$ws.Range("E23").Value = "'1.23%"
$ws.Range("D24").Value = "'0.001219"
$ws.Range("E24").Value = "'-0.04%"
$ws.Range("D25").Value = "'0.004288"
$ws.Range("E25").Value = "'-5.18%"
$ws.Range("E26").Value = "'0.09%"
$ws.Range("D27").Value = "'0.0001937"
$ws.Range("E27").Value = "'0.01%"
$ws.Range("D40").Value = "'0.03823"
$ws.Range("E40").Value = "'-0.38%"
$ws.Range("D41").Value = "'0.006270"
$ws.Range("E41").Value = "'56.21%"
$ws.Range("D42").Value = "'0.1103"
$ws.Range("E42").Value = "'-0.25%"
$ws.Range("D43").Value = "'0.002200"
$ws.Range("E43").Value = "'-6.69%"
$ws.Range("D44").Value = "'0.01144"
$ws.Range("E44").Value = "'6.88%"
$ws.Range("D45").Value = "'0.00005463"
$ws.Range("E45").Value = "'0.33%"
$ws.Range("E46").Value = "'0.06%"
$ws.Range("D47").Value = "'0.05000"
$ws.Range("E47").Value = "'-54.14%"
$ws.Range("D48").Value = "'0.2278"
$ws.Range("E48").Value = "'10,095.27%"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'0.06%"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'0.06%"
